# Sample Project rules sheet: the "R40" rule label in B11 is renamed to "1".
# The cell must keep holding a literal text value (not the number 1), so it
# is entered the same way a user typing into Excel would force text for a
# numeric-looking entry: with a leading apostrophe (quote prefix).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B11").Value = "'1"
